$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text values (e.g. '15.00', '64.811.86', '0.0000248') are not
# auto-converted to numbers/scientific notation by Excel's cell parser.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.852.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.147.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.73%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.64'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.01'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.46%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.146.56'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.03'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.61%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.56'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000248'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.676.57'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.31'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.475.59'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.154.19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '474.84'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.00'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.748'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.69'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.53'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +8.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '82.41'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.86'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.46'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.36%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.94%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +7.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.75'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0876'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.51'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.07'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.22'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.31%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '468.03'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.44'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.38'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.302'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.38%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.895.58'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.10'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.49'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.83'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.23%  '
